# Updates cryptos list cell values (Price / Volume(1h), and for rank-shifted rows also Coin/Link)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.833.04"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "2.372.98"
$ws.Range("E3").Value = "  -3.51%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "543.80"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "140.79"
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.556"
$ws.Range("E8").Value = "  -6.84%  "
$ws.Range("D9").Value = "2.371.97"
$ws.Range("E9").Value = "  -3.37%  "
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "5.34"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  -2.26%  "
$ws.Range("D14").Value = "25.38"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").Value = "2.796.47"
$ws.Range("E15").Value = "  -3.55%  "
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "60.704.78"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "2.368.36"
$ws.Range("E18").Value = "  -3.02%  "
$ws.Range("E19").Value = "  -4.07%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "316.43"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "6.67"
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  +3.66%  "
$ws.Range("D25").Value = "63.01"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.484.93"
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0925"
$ws.Range("E28").Value = "  -5.30%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "7.72"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "517.83"
$ws.Range("E30").Value = "  -2.49%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "7.95"
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "0.145"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "1.82"
$ws.Range("E34").Value = "  -2.94%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "5.45"
$ws.Range("E37").Value = "  -6.69%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "4.63"
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "0.374"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "18.02"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.72"
$ws.Range("E41").Value = "  +1.77%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "136.98"
$ws.Range("E43").Value = "  -5.12%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "40.18"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.23"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "139.00"
$ws.Range("E46").Value = "  -4.98%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "3.53"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "20.45"
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  -3.05%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.575"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.0915"
$ws.Range("E51").Value = "  -2.43%  "
